$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain alphabetic text can just be assigned directly.
function Set-PlainValue($rangeAddr, $text) {
    $ws.Range($rangeAddr).Value = $text
}

# Number-looking / date-looking text needs a leading apostrophe so Excel
# keeps it as literal text instead of silently auto-converting it into a
# real number or date (which would change the stored value, not just its
# display). Re-applying the "Normal" style afterwards keeps the cell look
# (general format, default alignment) as close as possible to the rest of
# the untouched sheet.
function Set-TextLookingValue($rangeAddr, $text) {
    $r = $ws.Range($rangeAddr)
    $r.Value = "'" + $text
    $r.Style = "Normal"
}

# Row 8
Set-PlainValue       "A8" "project 2"
Set-PlainValue       "B8" "edited"
Set-TextLookingValue "C8" "111"
Set-TextLookingValue "D8" "11/8/2021"
Set-TextLookingValue "E8" "12/8/2021"

# Row 9
Set-PlainValue       "A9" "test 4"
Set-TextLookingValue "B9" "1111"
Set-TextLookingValue "C9" "1111"
Set-TextLookingValue "E9" "11/9/2021"

# Row 10
Set-PlainValue       "A10" "project 5"
Set-PlainValue       "B10" "test5"
Set-TextLookingValue "C10" "1223"
Set-TextLookingValue "E10" "11/10/2021"

# Row 11
Set-PlainValue       "A11" "final"
Set-PlainValue       "B11" "test final"
Set-TextLookingValue "C11" "9998"
Set-TextLookingValue "D11" "11/8/2021"
Set-TextLookingValue "E11" "10/11/2021"
